$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ron", 75, 89, 78),
    @("Arthur", 90, 76, 59),
    @("Kate", 58, 89, 87),
    @("Gina", 66, 79, 90),
    @("Carrie", 74, 62, 65),
    @("Peter", 79, 88, 76),
    @("Kaylie", 61, 85, 78),
    @("Hanna", 79, 72, 76),
    @("Rose", 80, 64, 60),
    @("Jason", 85, 87, 82),
    @("Miguel", 87, 60, 90),
    @("Sarah", 70, 78, 79),
    @("Adam", 79, 76, 88),
    @("John", 77, 72, 82),
    @("Macy", 75, 70, 67),
    @("William", 82, 66, 65),
    @("Zara", 80, 62, 88),
    @("Kathy", 75, 65, 89),
    @("Carmen", 66, 74, 61),
    @("Brian", 88, 79, 86),
    @("Kumar", 78, 77, 74),
    @("Sam", 64, 78, 58),
    @("Lucy", 70, 67, 88),
    @("Andrew", 62, 60, 77)
)

$row = 7
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}
